# Rename "Interventions coverages" sheet to "Interventions cost and coverage"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Interventions coverages")
$ws.Name = "Interventions cost and coverage"

# Delete the four sheets related to birth distribution / time between births / RR birth by type / RR birth by time
$excel.DisplayAlerts = $false

$names = @("birth distribution", "time between births", "RR birth by type", "RR birth by time")
foreach ($name in $names) {
    $sheet = $wb.Worksheets.Item($name)
    $sheet.Delete()
}

$excel.DisplayAlerts = $true
